$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Add header row for new columns D, E, F
$ws.Range("D1").Value = "Duration"
$ws.Range("E1").Value = "Number of Modules"
$ws.Range("F1").Value = "Number of Topics"

$ws.Range("D2").Value = "4 місяці"
$ws.Range("E2").Value = 11
$ws.Range("F2").Value = 90
$ws.Range("D3").Value = "3 місяці"
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = 148
$ws.Range("D4").Value = "4 місяці"
$ws.Range("E4").Value = 19
$ws.Range("F4").Value = 185
$ws.Range("D5").Value = "5 місяців"
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 138
$ws.Range("D6").Value = "5 місяців"
$ws.Range("E6").Value = 25
$ws.Range("F6").Value = 182
$ws.Range("D7").Value = "4 місяці"
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 166
$ws.Range("D8").Value = "8 місяців"
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 90
$ws.Range("D9").Value = "5 місяців"
$ws.Range("E9").Value = 19
$ws.Range("F9").Value = 148
$ws.Range("D10").Value = "7 місяців"
$ws.Range("E10").Value = 19
$ws.Range("F10").Value = 185
$ws.Range("D11").Value = "5 місяців"
$ws.Range("E11").Value = 12
$ws.Range("F11").Value = 35
$ws.Range("D12").Value = "4 місяці"
$ws.Range("E12").Value = 9
$ws.Range("F12").Value = 64
$ws.Range("D13").Value = "7 місяців"
$ws.Range("E13").Value = 20
$ws.Range("F13").Value = 138
$ws.Range("D14").Value = "8 місяців"
$ws.Range("E14").Value = 25
$ws.Range("F14").Value = 182
$ws.Range("D15").Value = "4 місяці"
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 192
$ws.Range("D16").Value = "8 місяців"
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = 166
$ws.Range("D17").Value = "4 місяці"
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 45

# Copy formatting from existing columns so new cells match header/body styles
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1:F1").PasteSpecial(-4122)
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2:F17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set explicit column widths to match target layout (compensating for the
# Excel ColumnWidth -> stored width padding of ~0.8333 characters)
$ws.Columns.Item(4).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 18.166666666666666
$ws.Columns.Item(6).ColumnWidth = 17.166666666666666
